# Apply the two changes captured in the target diff:
#
#  1. The table on slide 16 switches from the default/themed table
#     style {AEDF82FB-2622-4CEB-B64B-23885AB6ECEF} to the built-in
#     style {F3720B07-0207-4359-863F-010E1052BFB8}.
#
#  2. The presentation's colour theme is swapped from the "Integral"
#     palette to the classic "Office" palette (the two theme parts in
#     the package effectively trade places).

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{F3720B07-0207-4359-863F-010E1052BFB8}")

# --- 2. Theme colours -------------------------------------------------
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink -> classic "Office" RGBs
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$themeColors = $slide.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
